$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_7_9_0",0.9999939966482478,0.9990518915967637,0.9999790192421294,0.999964122463273,0.9999714797316644,[double]"5.603865343999815e-06",0.0008850175772721206,[double]"1.015812931816948e-05",[double]"2.46499011003016e-05",[double]"1.740401520923553e-05",0.0001505861998536305,0.002367248475340055,1.000144080442052,0.002468027143815176,74.18410791665855,104.6560035383636),
    @("model_7_9_1",0.9999939501478314,0.9990514535360803,0.9999787096690488,0.9999638881510491,0.9999712250176203,[double]"5.64727144153437e-06",0.0008854264876914466,[double]"1.030801348369486e-05",[double]"2.481088687787095e-05",[double]"1.75594501807829e-05",0.0001498541312727336,0.002376398838901915,1.000145196452046,0.002477567057297728,74.1686761203128,104.6405717420178),
    @("model_7_9_2",0.9999939410612854,0.9990513706005084,0.9999786530568215,0.9999638408690785,0.9999711759421207,[double]"5.655753333435408e-06",0.0008855039044072871,[double]"1.033542309061084e-05",[double]"2.484337227139475e-05",[double]"1.75893976810028e-05",0.0001497201953009528,0.002378182779652441,1.000145414529151,0.00247942694409927,74.16567448361394,104.637570105319),
    @("model_7_9_3",0.9999939323352506,0.9990512892705907,0.9999785963631149,0.9999637969906949,0.9999711287503334,[double]"5.663898703880204e-06",0.0008855798223155132,[double]"1.036287214685952e-05",[double]"2.487351920773321e-05",[double]"1.761819567729637e-05",0.0001495906041641248,0.002379894683358951,1.000145623953986,0.002481211727090696,74.16279617186945,104.6346917935745),
    @("model_7_9_4",0.9999938905903151,0.9990509045263504,0.999978325560233,0.9999635904164091,0.9999709050315481,[double]"5.702865768811433e-06",0.0008859389641755696,[double]"1.049398517485078e-05",[double]"2.501544744976409e-05",[double]"1.775471631230744e-05",0.0001489652524601419,0.002388067371078847,1.000146625832438,0.002489732343046623,74.14908348602066,104.6209791077257),
    @("model_7_9_5",0.9999938866278237,0.9990508677167994,0.9999783031422204,0.9999635686374476,0.999970883877851,[double]"5.706564580545007e-06",0.0008859733242756131,[double]"1.050483917132163e-05",[double]"2.503041083063618e-05",[double]"1.77676250009789e-05",0.0001489004028368799,0.002388841681766501,1.00014672093223,0.002490539617743296,74.14778673000397,104.619682351709),
    @("model_7_9_6",0.9999938746230871,0.9990507618539561,0.9999782230615587,0.9999635119457332,0.9999708201951492,[double]"5.717770475260074e-06",0.0008860721425929415,[double]"1.054361135115791e-05",[double]"2.506936125142914e-05",[double]"1.780648630129352e-05",0.0001487330531663296,0.002391185997629643,1.00014700904591,0.002492983735986135,74.14386321121565,104.6157588329207),
    @("model_7_9_7",0.9999938668456743,0.9990506926952509,0.9999781746404706,0.9999634730680648,0.9999707791002763,[double]"5.725030348715882e-06",0.0008861366992085572,[double]"1.056705510270504e-05",[double]"2.509607241305042e-05",[double]"1.783156375787773e-05",0.0001486115898165751,0.002392703564739243,1.000147195703818,0.002494565909069387,74.14132541490456,104.6132210366096),
    @("model_7_9_8",0.9999938627419531,0.9990506586803101,0.9999781537736572,0.9999634495881189,0.9999707576043841,[double]"5.728860992230164e-06",0.0008861684506627604,[double]"1.057715806419591e-05",[double]"2.51122044666546e-05",[double]"1.784468126542526e-05",0.0001485698127416292,0.002393503915231843,1.000147294193125,0.002495400332139411,74.13998765342487,104.6118832751299),
    @("model_7_9_9",0.999993858849638,0.9990506249699688,0.9999781257731554,0.9999634325201829,0.9999707368881435,[double]"5.732494297615639e-06",0.0008861999177865971,[double]"1.059071490137145e-05",[double]"2.512393110594956e-05",[double]"1.785732300366051e-05",0.0001485003063943688,0.002394262787919413,1.000147387608688,0.002496191511608373,74.1387196339227,104.6106152556277),
    @("model_7_9_10",0.999993855570528,0.9990505922941875,0.9999781052675676,0.9999634159939464,0.9999707194501137,[double]"5.735555202806595e-06",0.0008862304191941461,[double]"1.06006429703783e-05",[double]"2.51352855670506e-05",[double]"1.786796426871445e-05",0.0001484491978537431,0.002394901919245671,1.000147466307328,0.002496857852078373,74.13765200502434,104.6095476267294),
    @("model_7_9_11",0.9999938513238465,0.9990505586160481,0.9999780778617221,0.9999633959713009,0.9999706973064338,[double]"5.739519293533747e-06",0.0008862618562589856,[double]"1.061391189634946e-05",[double]"2.514904225935535e-05",[double]"1.788147707785241e-05",0.0001483935794510625,0.002395729386540505,1.000147568227683,0.002497720546369043,74.13627019582177,104.6081658175268),
    @("model_7_9_12",0.9999938479900485,0.9990505261539965,0.9999780563150162,0.9999633797029278,0.9999706796005514,[double]"5.742631247639729e-06",0.0008862921581592883,[double]"1.06243440373668e-05",[double]"2.516021955362741e-05",[double]"1.78922817954971e-05",0.0001483355399831867,0.002396378777998113,1.000147648238835,0.002498397583765427,74.13518609431271,104.6070817160177),
    @("model_7_9_13",0.9999938446548937,0.9990504946330464,0.9999780357732676,0.9999633626085168,0.999970661828325,[double]"5.745744468246829e-06",0.0008863215815825863,[double]"1.063428961414666e-05",[double]"2.51719643828164e-05",[double]"1.790312699848153e-05",0.0001482818427345448,0.002397028257707203,1.000147728282551,0.002499074713170341,74.13410213927493,104.6059977609799),
    @("model_7_9_14",0.999993837064146,0.9990504315038885,0.9999779877524126,0.9999633260030443,0.9999706221713498,[double]"5.752830098022927e-06",0.0008863805098804107,[double]"1.065753958719104e-05",[double]"2.51971144171207e-05",[double]"1.792732700215587e-05",0.0001481810918696061,0.00239850580529273,1.000147910460496,0.002500615163015529,74.1316372665615,104.6035328882665),
    @("model_7_9_15",0.9999938337248501,0.9990503986240644,0.9999779671747381,0.999963309437757,0.9999706046827395,[double]"5.755947184215176e-06",0.0008864112017529296,[double]"1.066750255805707e-05",[double]"2.520849570831544e-05",[double]"1.793799913318626e-05",0.0001481227247387228,0.002399155514804152,1.000147990603598,0.002501292532005937,74.13055388952081,104.6024495112258),
    @("model_7_9_16",0.9999938301682546,0.9990503680458536,0.9999779445569016,0.9999632923107548,0.9999705860685733,[double]"5.759267109964193e-06",0.0008864397451705698,[double]"1.067845330202682e-05",[double]"2.522026292951388e-05",[double]"1.794935811577035e-05",0.0001480823599722712,0.002399847309718723,1.000148075961889,0.002502013778062234,74.12940065841067,104.6012962801157),
    @("model_7_9_17",0.9999938231396089,0.9990503065786638,0.9999778984564551,0.9999632596995413,0.9999705494219843,[double]"5.765828042263745e-06",0.0008864971220941247,[double]"1.070077348225472e-05",[double]"2.524266868143255e-05",[double]"1.797172108184363e-05",0.000147986975560263,0.002401213868497295,1.000148244649385,0.002503438514077176,74.1271235638022,104.5990191855072),
    @("model_7_9_18",0.9999938195370259,0.9990502751143922,0.9999778747822379,0.9999632428318855,0.9999705305347455,[double]"5.769190895327937e-06",0.0008865264926104847,[double]"1.071223568787941e-05",[double]"2.525425771691804e-05",[double]"1.798324670239873e-05",0.000147924910809789,0.002401914006647186,1.000148331111378,0.002504168458557594,74.12595742677975,104.5978530484848),
    @("model_7_9_19",0.9999938157646645,0.9990502461943167,0.9999778561785085,0.9999632213795157,0.9999705110780439,[double]"5.772712229114989e-06",0.0008865534882315364,[double]"1.072124294532996e-05",[double]"2.526899671078964e-05",[double]"1.79951198280598e-05",0.0001478812420148986,0.002402646921442056,1.000148421648052,0.002504932575052624,74.12473706169496,104.5966326834),
    @("model_7_9_20",0.9999938123718548,0.999050216827392,0.9999778355322693,0.9999632044992842,0.9999704933849372,[double]"5.7758792680985e-06",0.000886580900966689,[double]"1.073123911275261e-05",[double]"2.52805943864748e-05",[double]"1.80059167496137e-05",0.0001478339367515806,0.002403305903978622,1.000148503075485,0.002505619611840054,74.12364011785945,104.5955357395645),
    @("model_7_9_21",0.9999938089937925,0.999050186909235,0.9999778114348363,0.9999631905678056,0.999970475982716,[double]"5.779032540913365e-06",0.0008866088282530414,[double]"1.074290622418648e-05",[double]"2.529016610187633e-05",[double]"1.80165361630314e-05",0.0001478024582008017,0.002403961842649206,1.000148584148979,0.002506303475177825,74.12254853950895,104.594444161214),
    @("model_7_9_22",0.9999938061894703,0.9990501576044175,0.9999777917853361,0.9999631779824891,0.9999704611028248,[double]"5.781650252587409e-06",0.0008866361830137638,[double]"1.075241980641096e-05",[double]"2.529881292756406e-05",[double]"1.802561636698751e-05",0.0001477372757719538,0.002404506238833122,1.000148651452712,0.002506871047434346,74.12164281040435,104.5935384321094),
    @("model_7_9_23",0.9999938021126304,0.9990501282221562,0.9999777727534601,0.9999631543781102,0.9999704402648109,[double]"5.785455803647377e-06",0.0008866636100649237,[double]"1.076163435709427e-05",[double]"2.531503047367022e-05",[double]"1.803833241538224e-05",0.0001476981120983905,0.002405297445981968,1.000148749296869,0.002507695937909542,74.12032681963578,104.5922224413408),
    @("model_7_9_24",0.9999937993457926,0.9990500996293882,0.9999777520640157,0.9999631431107211,0.9999704257142914,[double]"5.788038525383177e-06",0.0008866903001588577,[double]"1.077165144289506e-05",[double]"2.532277180858722e-05",[double]"1.804721162574114e-05",0.0001476525045144519,0.002405834268062365,1.000148815700978,0.002508255613617277,74.11943418627709,104.5913298079821)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value2 = $rowVals[$c]
    }
}
